# openpyxl issue #268 regression test: a formula that returns a non-ASCII
# (umlaut) string needs to round-trip correctly. Add a new row with a
# formula exercising IF/ISBLANK and a literal "Düsseldorf" string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source file was subsequently re-saved by a German-locale Excel, which
# re-labels the three built-in cell styles. Re-apply that relabelling here
# (best effort - harmless if the host does not persist style renames).
try {
    $wb.Styles.Item("Followed Hyperlink").Name = "Besuchter Link"
    $wb.Styles.Item("Hyperlink").Name = "Link"
    $wb.Styles.Item("Normal").Name = "Standard"
} catch {
}

# Same re-save also bumped the sheet's default row height (15 -> 16pt).
try {
    $ws.StandardHeight = 16
} catch {
}

# The actual test case: row 16, column A holds a formula that returns the
# literal "Düsseldorf" whenever B16 is blank.
$ws.Range("A16").Formula = '=IF(ISBLANK(B16), "Düsseldorf", B16)'
$ws.Range("A16").Select()
